$d = $word.ActiveDocument

# Locate the target cell: row 4, column 1 of the first table -- this is
# the cell whose first paragraph reads "3. queries on multiple tables".
$t = $d.Tables.Item(1)
$cell = $t.Cell(4, 1)
$rng = $cell.Range

# Find the existing run of text and collapse the range to just after it.
$found = $rng.Find.Execute("queries on multiple tables", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)

# Insert a blank paragraph followed by a new paragraph holding the note
# text (still unformatted at this point).
$rng.InsertAfter("`r`rChưa chạy được")

# Re-acquire the cell range and color just the newly-added note text red
# (FF0000) without touching the pre-existing "queries on multiple
# tables" run.
$cell2 = $t.Cell(4, 1)
$rng2 = $cell2.Range
$rng2.Find.ClearFormatting()
$rng2.Find.Replacement.ClearFormatting()
$rng2.Find.Replacement.Font.Color = 255
$rng2.Find.Execute("Chưa chạy được", $true, $false, $false, $false, $false, $true, 1, $false, "Chưa chạy được", 2)
